$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/unstyled) used to reset cell style after forcing text format,
# so numeric-looking strings (e.g. "10.00") are stored as text without leaving a residual style on the cell.
$blankStyle = $ws.Range("D3").Style

$ws.Range("D2").Value = '67.369.46'
$ws.Range("E2").Value = '  -2.56%  '
$ws.Range("D3").Value = '3.765.57'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.81'
$ws.Range("D5").Style = $blankStyle
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.61'
$ws.Range("D6").Style = $blankStyle
$ws.Range("E6").Value = '  -3.55%  '
$ws.Range("D7").Value = '3.764.21'
$ws.Range("E7").Value = '  -1.49%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.516'
$ws.Range("D9").Style = $blankStyle
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.36'
$ws.Range("D11").Style = $blankStyle
$ws.Range("E11").Value = '  -2.81%  '
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("D13").Style = $blankStyle
$ws.Range("E13").Value = '  -4.91%  '
$ws.Range("E14").Value = '  -3.45%  '
$ws.Range("D15").Value = '4.398.70'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").Value = '3.794.60'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '67.380.08'
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.19'
$ws.Range("D21").Style = $blankStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '455.58'
$ws.Range("D22").Style = $blankStyle
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.696'
$ws.Range("D23").Style = $blankStyle
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000151'
$ws.Range("D24").Style = $blankStyle
$ws.Range("E24").Value = '  +2.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.08'
$ws.Range("D25").Style = $blankStyle
$ws.Range("E25").Value = '  -2.23%  '
$ws.Range("E26").Value = '  -3.65%  '
$ws.Range("E27").Value = '  -6.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("D29").Style = $blankStyle
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("E30").Value = '  -2.06%  '
$ws.Range("E31").Value = '  -4.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.59'
$ws.Range("D32").Style = $blankStyle
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.16'
$ws.Range("D33").Style = $blankStyle
$ws.Range("E33").Value = '  -4.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.13'
$ws.Range("D34").Style = $blankStyle
$ws.Range("E34").Value = '  -3.23%  '
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").Value = '3.720.04'
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0993'
$ws.Range("D37").Style = $blankStyle
$ws.Range("E37").Value = '  -3.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.29'
$ws.Range("D38").Style = $blankStyle
$ws.Range("E38").Value = '  -7.37%  '
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.992'
$ws.Range("D40").Style = $blankStyle
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.71'
$ws.Range("D41").Style = $blankStyle
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = $blankStyle
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.64'
$ws.Range("D44").Style = $blankStyle
$ws.Range("E44").Value = '  -1.45%  '
$ws.Range("E45").Value = '  -4.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.73'
$ws.Range("D46").Style = $blankStyle
$ws.Range("E46").Value = '  +0.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.32'
$ws.Range("D47").Style = $blankStyle
$ws.Range("E47").Value = '  -3.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '147.47'
$ws.Range("D48").Style = $blankStyle
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("E49").Value = '  -7.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '388.54'
$ws.Range("D50").Style = $blankStyle
$ws.Range("E50").Value = '  -4.06%  '
$ws.Range("D51").Value = '2.738.16'
$ws.Range("E51").Value = '  +1.20%  '
